$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Correlation")
$ws.Range("E2").Value = 123
